# Update cryptocurrency price/volume figures on Sheet1.
# These cells are stored as text (inline strings), so we force the
# NumberFormat to Text before writing to avoid Excel auto-converting
# values like "16.50" or "0.0788" into numbers and dropping trailing
# zeros / the leading "0.0" formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "37.020.09" }
    @{ Cell = "D3"; Value = "2.014.44" }
    @{ Cell = "E3"; Value = "  -1.62%  " }
    @{ Cell = "E4"; Value = "  +0.20%  " }
    @{ Cell = "D5"; Value = "226.34" }
    @{ Cell = "E5"; Value = "  -1.91%  " }
    @{ Cell = "D8"; Value = "54.68" }
    @{ Cell = "E8"; Value = "  -4.06%  " }
    @{ Cell = "D9"; Value = "0.379" }
    @{ Cell = "E9"; Value = "  -1.15%  " }
    @{ Cell = "D10"; Value = "0.0788" }
    @{ Cell = "E10"; Value = "  +2.65%  " }
    @{ Cell = "D11"; Value = "0.104" }
    @{ Cell = "E11"; Value = "  -3.34%  " }
    @{ Cell = "D12"; Value = "2.313.18" }
    @{ Cell = "E12"; Value = "  -1.54%  " }
    @{ Cell = "D13"; Value = "14.23" }
    @{ Cell = "E13"; Value = "  -2.95%  " }
    @{ Cell = "D14"; Value = "20.23" }
    @{ Cell = "E14"; Value = "  -1.65%  " }
    @{ Cell = "D15"; Value = "0.739" }
    @{ Cell = "E15"; Value = "  -2.08%  " }
    @{ Cell = "D16"; Value = "5.12" }
    @{ Cell = "E16"; Value = "  -2.29%  " }
    @{ Cell = "D17"; Value = "2.041.57" }
    @{ Cell = "E17"; Value = "  -0.59%  " }
    @{ Cell = "D18"; Value = "36.940.90" }
    @{ Cell = "D19"; Value = "6.07" }
    @{ Cell = "E19"; Value = "  +1.08%  " }
    @{ Cell = "D20"; Value = "68.79" }
    @{ Cell = "E20"; Value = "  -1.23%  " }
    @{ Cell = "D21"; Value = "0.0₃0818" }
    @{ Cell = "E21"; Value = "  -0.32%  " }
    @{ Cell = "D22"; Value = "223.56" }
    @{ Cell = "E22"; Value = "  -1.33%  " }
    @{ Cell = "E23"; Value = "  -0.06%  " }
    @{ Cell = "E24"; Value = "  +1.86%  " }
    @{ Cell = "D25"; Value = "2.19" }
    @{ Cell = "E25"; Value = "  -5.88%  " }
    @{ Cell = "D26"; Value = "164.87" }
    @{ Cell = "E26"; Value = "  -2.77%  " }
    @{ Cell = "E27"; Value = "  -4.87%  " }
    @{ Cell = "E28"; Value = "  -3.11%  " }
    @{ Cell = "D29"; Value = "1.35" }
    @{ Cell = "E29"; Value = "  +0.84%  " }
    @{ Cell = "D30"; Value = "18.70" }
    @{ Cell = "E30"; Value = "  -2.26%  " }
    @{ Cell = "E31"; Value = "  -3.47%  " }
    @{ Cell = "E32"; Value = "  -0.32%  " }
    @{ Cell = "D33"; Value = "0.0612" }
    @{ Cell = "E33"; Value = "  -1.54%  " }
    @{ Cell = "E34"; Value = "  -3.09%  " }
    @{ Cell = "D35"; Value = "2.34" }
    @{ Cell = "E35"; Value = "  -6.15%  " }
    @{ Cell = "E36"; Value = "  +2.01%  " }
    @{ Cell = "E37"; Value = "  +0.32%  " }
    @{ Cell = "E38"; Value = "  -4.59%  " }
    @{ Cell = "D39"; Value = "5.37" }
    @{ Cell = "E39"; Value = "  +1.95%  " }
    @{ Cell = "E40"; Value = "  -3.69%  " }
    @{ Cell = "D41"; Value = "1.472.99" }
    @{ Cell = "E41"; Value = "  -1.00%  " }
    @{ Cell = "D42"; Value = "94.98" }
    @{ Cell = "E42"; Value = "  -3.38%  " }
    @{ Cell = "D43"; Value = "16.50" }
    @{ Cell = "D44"; Value = "0.0920" }
    @{ Cell = "E44"; Value = "  -2.98%  " }
    @{ Cell = "E45"; Value = "  -4.97%  " }
    @{ Cell = "E46"; Value = "  -4.11%  " }
    @{ Cell = "E47"; Value = "  -0.13%  " }
    @{ Cell = "E48"; Value = "  -1.44%  " }
    @{ Cell = "D49"; Value = "2.91" }
    @{ Cell = "E49"; Value = "  -0.69%  " }
    @{ Cell = "D50"; Value = "2.202.85" }
    @{ Cell = "E50"; Value = "  -1.48%  " }
    @{ Cell = "D51"; Value = "44.17" }
    @{ Cell = "E51"; Value = "  -2.67%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.NumberFormat = "@"
    $range.Value = $u.Value
}
